$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1633333333333333
$ws.Range("C2").Value = 0.63
$ws.Range("J2").Value = 0.02666666666666667
$ws.Range("P2").Value = 0.1233333333333333
$ws.Range("S2").Value = 0.05666666666666666
$ws.Range("B3").Value = 0.005076142131979695
$ws.Range("C3").Value = 0.01522842639593909
$ws.Range("J3").Value = 0.04568527918781726
$ws.Range("P3").Value = 0.8071065989847716
$ws.Range("S3").Value = 0.1269035532994924
$ws.Range("J4").Value = 0.04166666666666666
$ws.Range("P4").Value = 0.75
$ws.Range("S4").Value = 0.2083333333333333
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.0423728813559322
$ws.Range("D6").Value = 0.02542372881355932
$ws.Range("F6").Value = 0.05084745762711865
$ws.Range("J6").Value = 0.2796610169491525
$ws.Range("O6").Value = 0.01694915254237288
$ws.Range("Q6").Value = 0.1991525423728814
$ws.Range("R6").Value = 0.09322033898305085
$ws.Range("S6").Value = 0.2923728813559322
$ws.Range("B7").Value = 0.07281553398058252
$ws.Range("D7").Value = 0.05825242718446602
$ws.Range("F7").Value = 0.07281553398058252
$ws.Range("J7").Value = 0.145631067961165
$ws.Range("O7").Value = 0.02427184466019417
$ws.Range("Q7").Value = 0.2233009708737864
$ws.Range("R7").Value = 0.07766990291262135
$ws.Range("S7").Value = 0.325242718446602
$ws.Range("B8").Value = 0.08713692946058091
$ws.Range("D8").Value = 0.02282157676348548
$ws.Range("F8").Value = 0.07883817427385892
$ws.Range("J8").Value = 0.1327800829875519
$ws.Range("O8").Value = 0.02282157676348548
$ws.Range("Q8").Value = 0.1701244813278008
$ws.Range("R8").Value = 0.1037344398340249
$ws.Range("S8").Value = 0.3817427385892116
$ws.Range("B9").Value = 0.08620689655172414
$ws.Range("D9").Value = 0.02873563218390805
$ws.Range("E9").Value = 0.005747126436781609
$ws.Range("F9").Value = 0.09195402298850575
$ws.Range("J9").Value = 0.1379310344827586
$ws.Range("O9").Value = 0.01724137931034483
$ws.Range("Q9").Value = 0.2586206896551724
$ws.Range("R9").Value = 0.05747126436781609
$ws.Range("S9").Value = 0.3160919540229885
$ws.Range("B10").Value = 0.1069057104913679
$ws.Range("D10").Value = 0.02523240371845949
$ws.Range("F10").Value = 0.06175298804780877
$ws.Range("J10").Value = 0.1248339973439575
$ws.Range("O10").Value = 0.01128818061088977
$ws.Range("Q10").Value = 0.2177954847277556
$ws.Range("R10").Value = 0.100265604249668
$ws.Range("S10").Value = 0.351925630810093
$ws.Range("G11").Value = 0.1438127090301003
$ws.Range("J11").Value = 0.08695652173913043
$ws.Range("K11").Value = 0.2073578595317726
$ws.Range("L11").Value = 0.5518394648829431
$ws.Range("S11").Value = 0.01003344481605351
$ws.Range("G12").Value = 0.7543859649122807
$ws.Range("J12").Value = 0.1637426900584795
$ws.Range("K12").Value = 0.005847953216374269
$ws.Range("L12").Value = 0.04678362573099415
$ws.Range("S12").Value = 0.02923976608187134
$ws.Range("F15").Value = 0.01901140684410646
$ws.Range("H15").Value = 0.1596958174904943
$ws.Range("I15").Value = 0.05703422053231939
$ws.Range("J15").Value = 0.4258555133079848
$ws.Range("K15").Value = 0.0532319391634981
$ws.Range("M15").Value = 0.007604562737642586
$ws.Range("O15").Value = 0.06083650190114068
$ws.Range("S15").Value = 0.2167300380228137
$ws.Range("F16").Value = 0.03017241379310345
$ws.Range("H16").Value = 0.1939655172413793
$ws.Range("I16").Value = 0.06465517241379311
$ws.Range("J16").Value = 0.4482758620689655
$ws.Range("K16").Value = 0.09051724137931035
$ws.Range("M16").Value = 0.02586206896551724
$ws.Range("N16").Value = 0.004310344827586207
$ws.Range("O16").Value = 0.05603448275862069
$ws.Range("S16").Value = 0.08620689655172414
$ws.Range("F17").Value = 0.01642335766423358
$ws.Range("H17").Value = 0.1751824817518248
$ws.Range("I17").Value = 0.08029197080291971
$ws.Range("J17").Value = 0.4416058394160584
$ws.Range("K17").Value = 0.1003649635036496
$ws.Range("M17").Value = 0.01642335766423358
$ws.Range("O17").Value = 0.06386861313868614
$ws.Range("S17").Value = 0.1058394160583942
$ws.Range("F18").Value = 0.0163265306122449
$ws.Range("H18").Value = 0.2122448979591837
$ws.Range("I18").Value = 0.07755102040816327
$ws.Range("J18").Value = 0.4693877551020408
$ws.Range("K18").Value = 0.07755102040816327
$ws.Range("M18").Value = 0.00816326530612245
$ws.Range("O18").Value = 0.05714285714285714
$ws.Range("S18").Value = 0.08163265306122448
$ws.Range("F19").Value = 0.01278976818545164
$ws.Range("H19").Value = 0.1998401278976819
$ws.Range("I19").Value = 0.06314948041566747
$ws.Range("J19").Value = 0.3956834532374101
$ws.Range("K19").Value = 0.09672262190247802
$ws.Range("M19").Value = 0.026378896882494
$ws.Range("N19").Value = 0.001598721023181455
$ws.Range("O19").Value = 0.08713029576338929
$ws.Range("S19").Value = 0.1167066346922462
